# Revert "Merging 0.1.8 w VitalSigns"
#
# - Rename the two "Include ValueSet #N" sheets to "Include ValueSets" / "Include ValueSets 2"
# - Roll back the Metadata sheet's Version/Status/Date/Contact values
# - Drop the "Jurisdiction" row (its contents were emptied by the merge)

$wb = $excel.ActiveWorkbook

# --- Rename the Include ValueSet sheets ---
$wb.Worksheets.Item(2).Name = "Include ValueSets"
$wb.Worksheets.Item(3).Name = "Include ValueSets 2"

$ws1 = $wb.Worksheets.Item(1)

# --- Roll back metadata values ---
$ws1.Range("B3").Value = "0.1.6"
$ws1.Range("B6").Value = "active"
$ws1.Range("B8").Value = "2023-05-05T10:50:04-05:00"

# Both Contact rows collapse to the single pre-merge placeholder text
$ws1.Range("B10").Value = "No display for ContactDetail"
$ws1.Range("B11").Value = "No display for ContactDetail"

# Remove the "Jurisdiction" row entirely (row 12), shifting later rows up
$ws1.Rows.Item(12).Delete()
